$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "19.973.14"
$ws.Range("E2").Value = "  -7.94%  "

$ws.Range("D3").Value = "1.410.96"
$ws.Range("E3").Value = "  -8.29%  "

$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").Value = "'1.001"
$ws.Range("E5").Value = "  -0.07%  "

$ws.Range("E6").Value = "  -5.54%  "

$ws.Range("D7").Value = "'0.3706"
$ws.Range("E7").Value = "  -5.13%  "

$ws.Range("D8").Value = "'0.3065"
$ws.Range("E8").Value = "  -3.12%  "

$ws.Range("D9").Value = "'39.36"
$ws.Range("E9").Value = "  -8.10%  "

$ws.Range("D10").Value = "'0.9976"
$ws.Range("E10").Value = "  -5.16%  "

$ws.Range("D11").Value = "'0.06558"
$ws.Range("E11").Value = "  -8.35%  "

$ws.Range("D12").Value = "'1.001"
$ws.Range("E12").Value = "  -0.01%  "

$ws.Range("D13").Value = "'5.392"
$ws.Range("E13").Value = "  -3.86%  "

$ws.Range("B14").Value = "Solana"
$ws.Range("C14").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D14").Value = "'16.99"
$ws.Range("E14").Value = "  -8.27%  "

$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "'6.171"
$ws.Range("E15").Value = "  -6.58%  "

$ws.Range("D16").Value = "1.409.48"
$ws.Range("E16").Value = "  -8.96%  "

$ws.Range("E17").Value = "  -8.00%  "

$ws.Range("D18").Value = "'0.05763"
$ws.Range("E18").Value = "  -12.26%  "

$ws.Range("D19").Value = "'73.55"

$ws.Range("D20").Value = "'1.000"
$ws.Range("E20").Value = "  -0.06%  "

$ws.Range("D21").Value = "'5.602"
$ws.Range("E21").Value = "  -8.74%  "

$ws.Range("D22").Value = "'14.43"

$ws.Range("D23").Value = "'10.83"
$ws.Range("E23").Value = "  -0.03%  "

$ws.Range("D24").Value = "'2.316"
$ws.Range("E24").Value = "  -2.62%  "

$ws.Range("D25").Value = "19.970.39"
$ws.Range("E25").Value = "  -7.96%  "

$ws.Range("D26").Value = "'2.263"
$ws.Range("E26").Value = "  -3.65%  "

$ws.Range("D27").Value = "'138.52"
$ws.Range("E27").Value = "  -5.95%  "

$ws.Range("D28").Value = "'16.84"
$ws.Range("E28").Value = "  -8.10%  "

$ws.Range("D29").Value = "1.566.92"
$ws.Range("E29").Value = "  -8.79%  "

$ws.Range("D30").Value = "'108.92"
$ws.Range("E30").Value = "  -6.98%  "

$ws.Range("E31").Value = "  -20.91%  "

$ws.Range("D32").Value = "'5.379"
$ws.Range("E32").Value = "  -8.10%  "

$ws.Range("D33").Value = "'0.8428"
$ws.Range("E33").Value = "  -12.11%  "

$ws.Range("D34").Value = "'0.07698"
$ws.Range("E34").Value = "  -5.74%  "

$ws.Range("D35").Value = "'8.430"
$ws.Range("E35").Value = "  -3.18%  "

$ws.Range("D36").Value = "'0.05794"
$ws.Range("E36").Value = "  -4.27%  "

$ws.Range("D37").Value = "'4.802"
$ws.Range("E37").Value = "  -5.65%  "

$ws.Range("D38").Value = "'0.9998"
$ws.Range("E38").Value = "  -0.08%  "

$ws.Range("D39").Value = "'0.1931"
$ws.Range("E39").Value = "  -4.71%  "

$ws.Range("D40").Value = "'0.02047"
$ws.Range("E40").Value = "  -6.64%  "

$ws.Range("D41").Value = "'10.25"
$ws.Range("E41").Value = "  -3.26%  "

$ws.Range("D42").Value = "'1.064"
$ws.Range("E42").Value = "  -9.32%  "

$ws.Range("D43").Value = "'1.281"
$ws.Range("E43").Value = "  -9.78%  "

$ws.Range("D44").Value = "'0.5293"
$ws.Range("E44").Value = "  -7.08%  "

$ws.Range("E45").Value = "  -5.31%  "

$ws.Range("D46").Value = "'12.15"
$ws.Range("E46").Value = "  -6.15%  "

$ws.Range("D47").Value = "'0.5119"
$ws.Range("E47").Value = "  -6.13%  "

$ws.Range("D48").Value = "'1.803"
$ws.Range("E48").Value = "  -2.74%  "

$ws.Range("D49").Value = "'110.06"
$ws.Range("E49").Value = "  -4.90%  "

$ws.Range("D50").Value = "'1.044"
$ws.Range("E50").Value = "  -9.89%  "

$ws.Range("D51").Value = "'1.001"
$ws.Range("E51").Value = "  +0.01%  "
